# Update the English back-cover translation paragraph:
#  - "This plugin aims" -> "This complement aims"
#  - "...to facilitate the first uses for technicians." ->
#    "...to facilitate the first uses of the tool " + "by the " + "technicians."
#    (split into 3 runs, matching the authored edit)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$oldText = "A very important aspect that must be considered in web development is the security of the interactions that originate between the client and the server. Being able to expand the functionalities that a tool offers us perfectly represents what could be a work environment where, after the solution has been delivered over time, the client asks us for new functionalities. Being able to provide the tool with security in its transactions and communication with the user is an ideal complementary aspect of the web development carried out in the final degree project. This would greatly benefit future users who are going to use the application. This plugin aims to cover essential fields in the security of the website, such as the tokenization of queries or the protection of access to certain routes of the page. The functionalities offered by the website will also be expanded, implementing an email manager that will help in all interactions with the user and a system capable of dumping data from an Excel file to facilitate the first uses for technicians. A backup system will be developed for the database that will allow the data stored in UAL Inventarium to be recovered with the minimum loss of information."

$newRun1 = "A very important aspect that must be considered in web development is the security of the interactions that originate between the client and the server. Being able to expand the functionalities that a tool offers us perfectly represents what could be a work environment where, after the solution has been delivered over time, the client asks us for new functionalities. Being able to provide the tool with security in its transactions and communication with the user is an ideal complementary aspect of the web development carried out in the final degree project. This would greatly benefit future users who are going to use the application. This complement aims to cover essential fields in the security of the website, such as the tokenization of queries or the protection of access to certain routes of the page. The functionalities offered by the website will also be expanded, implementing an email manager that will help in all interactions with the user and a system capable of dumping data from an Excel file to facilitate the first uses of the tool "
$newRun2 = "by the "
$newRun3 = "technicians. A backup system will be developed for the database that will allow the data stored in UAL Inventarium to be recovered with the minimum loss of information."

$fullText = $tr.Text
$startIdx0 = $fullText.IndexOf($oldText)
$start1 = $startIdx0 + 1

# Replace the whole old run's text with the full new text (single run for now).
$target = $tr.Characters($start1, $oldText.Length)
$newFull = $newRun1 + $newRun2 + $newRun3
$target.Text = $newFull

# Now split off the trailing run (run3) by re-setting its own sub-range text.
$run3Start = $start1 + $newRun1.Length + $newRun2.Length
$run3Range = $tr.Characters($run3Start, $newRun3.Length)
$run3Range.Text = $newRun3

# Then split off the middle run (run2) the same way.
$run2Start = $start1 + $newRun1.Length
$run2Range = $tr.Characters($run2Start, $newRun2.Length)
$run2Range.Text = $newRun2
